$wb = $excel.ActiveWorkbook

# The workbook's last test-case sheet ("tc_051") is the template for the
# new one: same layout (header row + value row in column A), same
# column width/styles. Duplicate it (preserves styles/col width/dimension)
# and rename/re-seed the copy, matching how the previous tc_0xx sheets
# were produced.
$src = $wb.Worksheets.Item("tc_051")
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item($src.Index + 1)
$newSheet.Name = "tc_053"

# New content for the duplicated sheet.
$newSheet.Range("A1").Value = "reject_comment_v7_user"
$newSheet.Range("A2").Value = "Reject comment as V7 user automation test"

# The previously-active sheet (tc_051) is no longer the selected tab;
# its selection collapses to the full used range (no explicit
# ActiveCell) the way the other superseded tc_0xx sheets look.
$src.Activate()
$src.Range("A1:A2").Select()

# The newly added sheet becomes the active / selected tab.
$newSheet.Select()
